$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date that was bumped from 2023-09-23 (45192)
# to 2023-10-03 (45202) for every data row (rows 2 through 261).
$oldValue = 45192
$newValue = 45202

for ($row = 2; $row -le 261; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
